$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row = 4; Place = "HaUI"; Height = 175.5185207047189; Timestamp = "2025-05-09 04:17:38"; Path = "/home/anodi108/Desktop/project/Do_An_Tot_Nghiep/DATN_PhamDangDong/DATN_PhamDangDong/resource/data/data_result/image_20250509_041738.jpg" },
    @{ Row = 5; Place = "HaUI"; Height = 175.5185207047189; Timestamp = "2025-05-09 04:18:38"; Path = "/home/anodi108/Desktop/project/Do_An_Tot_Nghiep/DATN_PhamDangDong/DATN_PhamDangDong/resource/data/data_result/image_20250509_041838.jpg" },
    @{ Row = 6; Place = "HaUI"; Height = 175.5185207047189; Timestamp = "2025-05-09 04:20:59"; Path = "/home/anodi108/Desktop/project/Do_An_Tot_Nghiep/DATN_PhamDangDong/DATN_PhamDangDong/resource/data/data_result/image_20250509_042059.jpg" }
)

foreach ($r in $rows) {
    $ws.Cells.Item($r.Row, 1).Value = $r.Place
    $ws.Cells.Item($r.Row, 2).Value = $r.Height
    $ws.Cells.Item($r.Row, 3).Value = $r.Timestamp
    $ws.Cells.Item($r.Row, 4).Value = $r.Path
}
